# Actualizar precios con datos nuevos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: tiny recompute of the timestamp (same moment, refreshed precision)
$ws.Cells.Item(7, 1).Value = 45807.39290912037

# New row 8: same product/weight/price pattern as the previous entries,
# timestamped the next day.
$ws.Cells.Item(8, 1).Value = 45808.3913043522
$ws.Cells.Item(8, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(8, 2).Value = "EVOWHEY PROTEIN"
$ws.Cells.Item(8, 3).Value = "2Kg"
$ws.Cells.Item(8, 4).Value = "37,90€"
